# Auto - Update data with bot!
# Applies the row-level content corrections as described in the commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 26: title update
$ws.Range("D26").Value = "ai plus(est soft)"

# Row 29: title + link update
$ws.Range("D29").Value = "[만화] 인턴일기 66~71"
$ws.Range("E29").Value = "https://blog.promedius.ai/intern-life-10/"

# Row 37: title + link update
$ws.Range("D37").Value = "[Paper Review] Anomaly Transformer: Time Series Anomaly Detection with Association Discrepancy"
$ws.Range("E37").Value = "http://dsba.korea.ac.kr/seminar/?uid=1936&mod=document&pageid=1"

# Row 51: title + link update
$ws.Range("D51").Value = "[jquery] 체크박스 체크 여부 확인하기"
$ws.Range("E51").Value = "https://bskyvision.com/1246"

# Row 52: title + link update
$ws.Range("D52").Value = "[R] for: 결과가 가변적일 때"
$ws.Range("E52").Value = "http://ds.sumeun.org/?p=2555&utm_source=rss&utm_medium=rss&utm_campaign=for-%25ea%25b2%25b0%25ea%25b3%25bc%25ea%25b0%2580-%25eb%25b3%2580%25ec%25a0%2581%25ec%259d%25bc-%25eb%2595%258c"

$wb.Save()
